$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values are plain decimals (e.g. "12.00",
# "1.43") that Excel would otherwise auto-convert to numbers (dropping
# trailing zeros / changing type) when assigned through .Value. Force
# those specific cells to Text format first so the literal string is
# preserved exactly, then restore the original (Normal) cell style so
# no unintended formatting change is left behind.
$textCells = @('D5', 'D6', 'D11', 'D12', 'D13', 'D14', 'D18', 'D19', 'D23', 'D25', 'D26', 'D28', 'D34', 'D36', 'D39', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.891.46'
$ws.Range('E2').Value = '  +1.67%  '
$ws.Range('D3').Value = '3.817.04'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '626.09'
$ws.Range('E5').Value = '  +4.32%  '
$ws.Range('D6').Value = '165.04'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = '3.815.99'
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('E10').Value = '  +1.23%  '
$ws.Range('D11').Value = '0.453'
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').Value = '6.63'
$ws.Range('E12').Value = '  +3.79%  '
$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').Value = '35.85'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').Value = '4.457.25'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('D16').Value = '3.839.12'
$ws.Range('E16').Value = '  +1.51%  '
$ws.Range('D17').Value = '68.891.44'
$ws.Range('D18').Value = '18.15'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('D19').Value = '7.13'
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('D23').Value = '0.708'
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('E24').Value = '  +4.45%  '
$ws.Range('D25').Value = '83.83'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = '12.00'
$ws.Range('E26').Value = '  -0.30%  '
$ws.Range('E27').Value = '  +2.06%  '
$ws.Range('D28').Value = '10.03'
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '3.968.98'
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('E31').Value = '  -3.75%  '
$ws.Range('E32').Value = '  +2.06%  '
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('D34').Value = '29.14'
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').Value = '9.10'
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('E37').Value = '  +1.73%  '
$ws.Range('E38').Value = '  +7.52%  '
$ws.Range('D39').Value = '3.35'
$ws.Range('E39').Value = '  +3.79%  '
$ws.Range('E40').Value = '  +2.60%  '
$ws.Range('D41').Value = '0.980'
$ws.Range('E41').Value = '  -1.74%  '
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D44').Value = '1.43'
$ws.Range('E44').Value = '  +3.03%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '0.300'
$ws.Range('E45').Value = '  +0.97%  '
$ws.Range('D46').Value = '154.63'
$ws.Range('E46').Value = '  +2.34%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '46.78'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('B48').Value = 'Arweave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').Value = '42.79'
$ws.Range('E48').Value = '  -5.52%  '
$ws.Range('D49').Value = '8.44'
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('D51').Value = '380.78'
$ws.Range('E51').Value = '  -3.47%  '

# Restore original (General/Normal) style on the cells we had to
# temporarily mark as Text, so formatting matches the source workbook.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
